$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Periods (column E) are now listed newest-first (descending) instead of
# oldest-first (ascending): the whole "Periodo Mora" table got reversed.
$periods = @(
    "2304","2303","2302","2301",
    "2212","2211","2210","2209","2208","2207","2206","2205","2204","2203","2202","2201",
    "2112","2111","2110","2109","2108","2107","2106","2105","2104","2103","2102","2101",
    "2012","2011","2010","2009","2008","2007","2006","2005","2004","2003","2002","2001",
    "1912","1911","1910","1909","1908","1907","1906","1905"
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}

# Column F ("Valor Mora") values travel with their period: period 2304 always
# carried 29812 and period 1905 always carried 26500, so after the reorder the
# two figures swap between the first and last data row.
$ws.Cells.Item(16, 6).Value = 29812
$ws.Cells.Item(63, 6).Value = 26500
